$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) and E (Volume) to be treated as text so that
# numeric-looking values (e.g. "27.642.09", "1.003") are not auto-converted
# into numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.642.09'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.842.95'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = '312.66'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.4280'
$ws.Range("E7").Value = '  +1.41%  '
$ws.Range("D8").Value = '0.3629'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.07317'
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("D10").Value = '0.8762'
$ws.Range("E10").Value = '  -1.26%  '
$ws.Range("D11").Value = '20.63'
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").Value = '1.871.66'
$ws.Range("E12").Value = '  +2.12%  '
$ws.Range("D13").Value = '5.344'
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("D14").Value = '6.510'
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").Value = '0.06962'
$ws.Range("E15").Value = '  +1.59%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '79.45'
$ws.Range("E17").Value = '  +1.36%  '
$ws.Range("D18").Value = '0.000008956'
$ws.Range("E18").Value = '  +1.84%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '15.34'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").Value = '27.639.21'
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("D22").Value = '4.979'
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("E23").Value = '  -1.94%  '
$ws.Range("D24").Value = '2.039.57'
$ws.Range("E24").Value = '  -0.84%  '
$ws.Range("D25").Value = '1.988'
$ws.Range("E25").Value = '  -2.22%  '
$ws.Range("D26").Value = '155.42'
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("D27").Value = '18.52'
$ws.Range("E27").Value = '  +0.53%  '
$ws.Range("D28").Value = '119.55'
$ws.Range("E28").Value = '  +0.91%  '
$ws.Range("D29").Value = '5.200'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '1.865'
$ws.Range("E30").Value = '  +3.47%  '
$ws.Range("D31").Value = '0.08892'
$ws.Range("E31").Value = '  +0.34%  '
$ws.Range("D32").Value = '0.7524'
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("D33").Value = '2.962'
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").Value = '4.517'
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("D35").Value = '1.124'
$ws.Range("E35").Value = '  +2.22%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '0.05436'
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("D38").Value = '1.105'
$ws.Range("E38").Value = '  +0.87%  '
$ws.Range("D39").Value = '0.01933'
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("D40").Value = '2.822'
$ws.Range("E40").Value = '  +1.52%  '
$ws.Range("D41").Value = '0.1663'
$ws.Range("E41").Value = '  +1.14%  '
$ws.Range("D42").Value = '0.5066'
$ws.Range("E42").Value = '  +0.52%  '
$ws.Range("D43").Value = '6.587'
$ws.Range("E43").Value = '  -3.22%  '
$ws.Range("D44").Value = '8.366'
$ws.Range("E44").Value = '  +2.80%  '
$ws.Range("D45").Value = '0.06547'
$ws.Range("E45").Value = '  -0.97%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '10.38'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '105.90'
$ws.Range("E47").Value = '  +0.88%  '
$ws.Range("D48").Value = '0.4649'
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("D49").Value = '1.000'
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").Value = '64.44'
$ws.Range("E51").Value = '  +0.15%  '
